$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) figures on the "展览" sheet
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 309
$wsExhibition.Range("F5").Value = 123

# Update the same rows on the "全部类型" sheet (duplicated data)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 309
$wsAll.Range("F5").Value = 123
